# Updated symbol list on Tue Jan 24 17:52:12 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the crypto ticker rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.77%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.82%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.101"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.47%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08164"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.061"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.09%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.941"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.23%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.129"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.35%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "10.83%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9253"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.30%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1070"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "12.22%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1915"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.73%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09183"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.19%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03656"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.19%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09911"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.11%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001427"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.42%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005657"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.90%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.474"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3394"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.74%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.099"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.68%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2215"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.55%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04538"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.74%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001227"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.54%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004780"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.08%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001251"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004453"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.35%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01959"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04883"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.35%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007561"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.10%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009989"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "29.18%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1382"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.17%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002172"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.90%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.98%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006592"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.51%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "182.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "253.93%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001502"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-21.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
